$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("H2").Value = 0.0321
$ws.Range("I2").Value = 0.0053
$ws.Range("J2").Value = 0.9899
$ws.Range("K2").Value = 0.0026
$ws.Range("L2").Value = 0.9899
$ws.Range("M2").Value = 0.0026
$ws.Range("N2").Value = 0.9916
$ws.Range("O2").Value = 0.0028
$ws.Range("P2").Value = 0.9882
$ws.Range("Q2").Value = 0.0048
$ws.Range("R2").Value = 0.9899
$ws.Range("S2").Value = 0.0026

# Row 3 updates
$ws.Range("H3").Value = 0.0324
$ws.Range("I3").Value = 0.0052
$ws.Range("J3").Value = 0.9898
$ws.Range("K3").Value = 0.0027
$ws.Range("L3").Value = 0.9898
$ws.Range("M3").Value = 0.0027
$ws.Range("N3").Value = 0.9913999999999999
$ws.Range("O3").Value = 0.0031
$ws.Range("P3").Value = 0.9882
$ws.Range("Q3").Value = 0.0046
$ws.Range("R3").Value = 0.9898
$ws.Range("S3").Value = 0.0027
